$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "Conception du menu" -> Terminée (completed)
$ws.Range("C7").Value2 = "Terminée"
$ws.Range("F7").Value2 = $ws.Range("E7").Value2
$ws.Range("G7").Value2 = 1
$ws.Range("H7").Value2 = "Non"

# Row 8: "Recherche mathématiques sur la détection des collisions" -> Terminée
$ws.Range("C8").Value2 = "Terminée"
$ws.Range("F8").Value2 = $ws.Range("E8").Value2
$ws.Range("G8").Value2 = 1
$ws.Range("H8").Value2 = "Non"

# Row 9: "Création du menu" -> Terminée
$ws.Range("C9").Value2 = "Terminée"
$ws.Range("F9").Value2 = $ws.Range("E9").Value2
$ws.Range("G9").Value2 = 1
$ws.Range("H9").Value2 = "Non"

# Row 10: "Conception du joueurs et des bunkers" -> Terminée
$ws.Range("C10").Value2 = "Terminée"
$ws.Range("F10").Value2 = $ws.Range("E10").Value2
$ws.Range("G10").Value2 = 1
$ws.Range("H10").Value2 = "Non"

# Row 11: "Conception du système de collisions" -> Terminée
$ws.Range("C11").Value2 = "Terminée"
$ws.Range("F11").Value2 = $ws.Range("E11").Value2
$ws.Range("G11").Value2 = 1
$ws.Range("H11").Value2 = "Non"

# Row 12: "Programmation du joueur (déplacements)" -> Terminée
$ws.Range("C12").Value2 = "Terminée"
$ws.Range("F12").Value2 = $ws.Range("E12").Value2
$ws.Range("G12").Value2 = 1
$ws.Range("H12").Value2 = "Non"

# Row 13: "Programmation des bunkers" -> Terminée
$ws.Range("C13").Value2 = "Terminée"
$ws.Range("F13").Value2 = $ws.Range("E13").Value2
$ws.Range("G13").Value2 = 1
$ws.Range("H13").Value2 = "Non"

# Row 14: "Conceptions des ennemies" -> En cours de réalisation
$ws.Range("C14").Value2 = "En cours de réalisation"

# Update selection to C14
$ws.Range("C14").Select()
